# Restore the "adductName" column (as column C) to the "Corrected" sheet.
#
# The sheet previously had columns: Compound, C_Label, blank_1_404020 (/
# formula), <sample1>, <sample2> (A:E). The edit re-inserts an "adductName"
# column between C_Label and the formula/blank column, shifting the
# remaining columns one slot to the right (A:F), and back-fills the new
# column's data rows with 0 (mirroring how this column is populated
# elsewhere in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

# Shift existing column C (and right) over to make room; the new column
# inherits column C's old formatting (bold+centered header style), matching
# how Excel itself behaves on a column insert.
$ws.Columns("C").Insert()

# Header text for the newly inserted column.
$ws.Range("C1").Value = "adductName"

# The restored column's header is bold but left-aligned (not centered like
# the other headers), so pull in a plain (non-centered) format and layer
# Bold back on top of it, rather than keeping the inherited centered style.
$ws.Range("A2").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Font.Bold = $true

# Data rows: the adduct name isn't populated by this loader (yet), so every
# row gets a placeholder 0, consistent with the other restored test files.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Match the column's original best-fit width.
$ws.Columns("C").ColumnWidth = 10.2

Write-Host "Inserted adductName column into Corrected sheet"
